$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings from "DD/MM/YYYY" to "DD-MM-YYYY" format.
# Force text interpretation (Excel auto-detects some of these as dates and
# would otherwise convert them to date serial numbers), then restore the
# original "Normal" style so no stray formatting is introduced.
$dates = [ordered]@{
  "A3"  = "28-07-2022"
  "A4"  = "01-08-2022"
  "A5"  = "04-08-2022"
  "A6"  = "08-08-2022"
  "A7"  = "11-08-2022"
  "A8"  = "15-08-2022"
  "A9"  = "18-08-2022"
  "A10" = "22-08-2022"
  "A11" = "25-08-2022"
  "A12" = "29-08-2022"
  "A13" = "01-09-2022"
  "A14" = "05-09-2022"
  "A15" = "08-09-2022"
  "A16" = "12-09-2022"
  "A17" = "15-09-2022"
  "A18" = "19-09-2022"
  "A19" = "22-09-2022"
  "A20" = "26-09-2022"
  "A21" = "29-09-2022"
}

foreach ($addr in $dates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$addr]
    $cell.Style = "Normal"
}

# Update attendance counts that changed
# Row 3: Total=1, Real=0, Duplicate=0, Invalid=1, Absent=1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: Total=1, Real=1, Duplicate=0, Invalid=0, Absent=0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 10: Total=1, Real=1, Duplicate=0, Invalid=0, Absent=0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0
